# Update cryptos list with latest prices/volumes (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.211.24"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.140.34"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.29"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.05"
$ws.Range("E6").Value = "  -4.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -7.50%  "
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.57"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.380"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "3.686.37"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "64.287.43"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.91"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "3.140.48"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000154"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "402.70"
$ws.Range("E18").Value = "  -4.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.53"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.22"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.11"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.87"
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.02"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.481"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.194"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000101"
$ws.Range("E27").Value = "  -4.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.78"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.15"
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.24"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").Value = "2.674.33"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.67"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.52"
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.07"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.688"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0612"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.48"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "287.45"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.08"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.996"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0974"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.46"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.89"
$ws.Range("E50").Value = "  -5.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.65"
$ws.Range("E51").Value = "  -1.93%  "

# Rows 33 and 34 swapped ranking order (NEARProtocol moved above Monero),
# with updated price/volume values
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.81"
$ws.Range("E33").Value = "  -3.87%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.03"
$ws.Range("E34").Value = "  +0.53%  "

